$wb = $excel.ActiveWorkbook

# --- Update "Not Scottish" sheet with computed comparison values ---
$ns = $wb.Worksheets.Item("Not Scottish")
$ns.Cells.Item(2, 3).Value = -0.0630372492836676
$ns.Cells.Item(2, 4).Value = 0.041751944330740955
$ns.Cells.Item(2, 5).Value = -0.02828989192625545
$ns.Cells.Item(2, 6).Value = 0.06989779096604021
$ns.Cells.Item(3, 3).Value = -0.0725933719095212
$ns.Cells.Item(3, 4).Value = 0.033609352167559614
$ns.Cells.Item(3, 5).Value = -0.047235023041474776
$ns.Cells.Item(3, 6).Value = 0.03525523319867792
$ns.Cells.Item(4, 3).Value = -0.06706204379562052
$ns.Cells.Item(4, 4).Value = 0.03883495145631069
$ns.Cells.Item(4, 5).Value = -0.03721571330117172
$ns.Cells.Item(4, 6).Value = 0.05463347164591971
$ns.Cells.Item(5, 3).Value = -0.0701606086221471
$ns.Cells.Item(5, 4).Value = 0.038299663299663285
$ns.Cells.Item(5, 5).Value = -0.02997099581050606
$ns.Cells.Item(5, 6).Value = 0.06423553027768483
$ns.Cells.Item(6, 3).Value = -0.04267161410018559
$ns.Cells.Item(6, 4).Value = 0.05987055016181233
$ns.Cells.Item(6, 5).Value = -0.013771186440677919
$ns.Cells.Item(6, 6).Value = 0.21515892420537902
$ns.Cells.Item(7, 3).Value = -0.05942376950780307
$ns.Cells.Item(7, 4).Value = 0.019499417927823
$ns.Cells.Item(7, 5).Value = -0.022898961284230374
$ns.Cells.Item(7, 6).Value = 0.06167400881057269
$ns.Cells.Item(8, 3).Value = -0.15300546448087435
$ns.Cells.Item(8, 4).Value = 0.025773195876288683
$ns.Cells.Item(8, 5).Value = -0.03793103448275863
$ns.Cells.Item(8, 6).Value = 0.12546125461254615
$ns.Cells.Item(9, 3).Value = -0.1152073732718895
$ns.Cells.Item(9, 4).Value = 0.021298174442190718
$ns.Cells.Item(9, 5).Value = -0.04664310954063601
$ns.Cells.Item(9, 6).Value = 0.12720588235294103
$ns.Cells.Item(10, 3).Value = -0.08175675675675675
$ns.Cells.Item(10, 4).Value = 0.011647254575707259
$ns.Cells.Item(10, 5).Value = -0.05378973105134481
$ns.Cells.Item(10, 6).Value = 0.12186084808563201
$ns.Cells.Item(11, 3).Value = 0.3161693936477382
$ns.Cells.Item(11, 4).Value = 0.24127906976744193
$ns.Cells.Item(11, 5).Value = 0.22485422740524796
$ns.Cells.Item(11, 6).Value = 0.45116279069767456
$ns.Cells.Item(12, 3).Value = 0.16666666666666663
$ns.Cells.Item(12, 4).Value = 0.17721518987341778
$ns.Cells.Item(12, 5).Value = 0.08316633266533068
$ns.Cells.Item(12, 6).Value = 0.4921383647798741
$ns.Cells.Item(13, 3).Value = 0.036728563626433536
$ns.Cells.Item(13, 4).Value = -0.02260317460317451
$ns.Cells.Item(13, 5).Value = 0.0012208521548040543
$ns.Cells.Item(13, 6).Value = -0.002864919066036387
$ns.Cells.Item(14, 3).Value = -0.013667117726657643
$ns.Cells.Item(14, 4).Value = -0.007851797325481625
$ns.Cells.Item(14, 5).Value = -0.007304983042003573
$ns.Cells.Item(14, 6).Value = -0.019829059829059796
$ns.Cells.Item(15, 3).Value = 0.5373711340206185
$ns.Cells.Item(15, 4).Value = 0.6246458923512749
$ns.Cells.Item(15, 5).Value = 0.29635145197319435
$ns.Cells.Item(15, 6).Value = 0.6837146702557197
$ns.Cells.Item(16, 3).Value = -0.09566849552411205
$ns.Cells.Item(16, 4).Value = -0.05856860434541258
$ns.Cells.Item(16, 5).Value = -0.06714445688689816
$ns.Cells.Item(16, 6).Value = -0.09791099182823271
$ns.Cells.Item(17, 3).Value = 0.41064638783269963
$ns.Cells.Item(17, 4).Value = 0.19811320754716985
$ns.Cells.Item(17, 5).Value = 0.019880715705765425
$ns.Cells.Item(17, 6).Value = 0.12451361867704272
$ns.Cells.Item(18, 3).Value = 0.3658675799086759
$ns.Cells.Item(18, 4).Value = 0.2915407854984896
$ns.Cells.Item(18, 5).Value = 0.2792718245759205
$ns.Cells.Item(18, 6).Value = 0.4772162386081192
$ns.Cells.Item(19, 3).Value = 0.0
$ns.Cells.Item(19, 4).Value = 0.09844559585492221
$ns.Cells.Item(19, 5).Value = 0.14814814814814806
$ns.Cells.Item(19, 6).Value = 1.9411764705882353
$ns.Cells.Item(20, 3).Value = 0.0484173819742488
$ns.Cells.Item(20, 4).Value = 0.08197029516826465
$ns.Cells.Item(20, 5).Value = 0.0261437908496732
$ns.Cells.Item(20, 6).Value = 0.11039603960396052

# --- Add new "Student" sheet after "Single" with full seg-comparison table ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$student = $wb.Worksheets.Add($null, $lastSheet)
$student.Name = "Student"

# Header row
$student.Cells.Item(1, 2).Value = "index"
$student.Cells.Item(1, 3).Value = "Aberdeen"
$student.Cells.Item(1, 4).Value = "Dundee"
$student.Cells.Item(1, 5).Value = "Edinburgh"
$student.Cells.Item(1, 6).Value = "Glasgow"

# Data rows (A: seg index as text, B: metric label, C-F: city values)
$student.Cells.Item(2, 1).Value = "'1"
$student.Cells.Item(2, 2).Value = "IS"
$student.Cells.Item(2, 3).Value = -0.015846538782318703
$student.Cells.Item(2, 4).Value = -0.03748594277146074
$student.Cells.Item(2, 5).Value = -0.06652719665271963
$student.Cells.Item(2, 6).Value = -0.09214501510574023
$student.Cells.Item(3, 1).Value = "'2"
$student.Cells.Item(3, 2).Value = "IS(adj)"
$student.Cells.Item(3, 3).Value = -0.019643608811561686
$student.Cells.Item(3, 4).Value = -0.04000000000000009
$student.Cells.Item(3, 5).Value = -0.07076006806579702
$student.Cells.Item(3, 6).Value = -0.09690904484418547
$student.Cells.Item(4, 1).Value = "'3"
$student.Cells.Item(4, 2).Value = "IS(w)"
$student.Cells.Item(4, 3).Value = -0.017311182465447546
$student.Cells.Item(4, 4).Value = -0.03845670478823673
$student.Cells.Item(4, 5).Value = -0.06856821694534207
$student.Cells.Item(4, 6).Value = -0.09429437010855843
$student.Cells.Item(5, 1).Value = "'4"
$student.Cells.Item(5, 2).Value = "IS(s)"
$student.Cells.Item(5, 3).Value = -0.016441410059913657
$student.Cells.Item(5, 4).Value = -0.0380904648540283
$student.Cells.Item(5, 5).Value = -0.0668064290705801
$student.Cells.Item(5, 6).Value = -0.092891353667759
$student.Cells.Item(6, 1).Value = "'5"
$student.Cells.Item(6, 2).Value = "H"
$student.Cells.Item(6, 3).Value = 0.010000000000000047
$student.Cells.Item(6, 4).Value = -0.03558627752176145
$student.Cells.Item(6, 5).Value = -0.0793795620437955
$student.Cells.Item(6, 6).Value = -0.1138034304383338
$student.Cells.Item(7, 1).Value = "'6"
$student.Cells.Item(7, 2).Value = "G"
$student.Cells.Item(7, 3).Value = -0.01492537313432841
$student.Cells.Item(7, 4).Value = -0.02247314185485635
$student.Cells.Item(7, 5).Value = -0.0417246175243394
$student.Cells.Item(7, 6).Value = -0.054319371727748644
$student.Cells.Item(8, 1).Value = "'7"
$student.Cells.Item(8, 2).Value = "A(0.1)"
$student.Cells.Item(8, 3).Value = -0.05446082234290134
$student.Cells.Item(8, 4).Value = -0.1247392574050897
$student.Cells.Item(8, 5).Value = -0.17879558948261245
$student.Cells.Item(8, 6).Value = -0.2241743725231175
$student.Cells.Item(9, 1).Value = "'8"
$student.Cells.Item(9, 2).Value = "A(0.5)"
$student.Cells.Item(9, 3).Value = -0.03050524308865583
$student.Cells.Item(9, 4).Value = -0.057486631016042795
$student.Cells.Item(9, 5).Value = -0.09642857142857134
$student.Cells.Item(9, 6).Value = -0.12491035142242414
$student.Cells.Item(10, 1).Value = "'9"
$student.Cells.Item(10, 2).Value = "A(0.9)"
$student.Cells.Item(10, 3).Value = -0.023980815347721843
$student.Cells.Item(10, 4).Value = -0.03372066529961275
$student.Cells.Item(10, 5).Value = -0.05892504930966477
$student.Cells.Item(10, 6).Value = -0.07622913847541722
$student.Cells.Item(11, 1).Value = "'10"
$student.Cells.Item(11, 2).Value = "xPx"
$student.Cells.Item(11, 3).Value = 0.275
$student.Cells.Item(11, 4).Value = 0.15840893230983946
$student.Cells.Item(11, 5).Value = 0.09799382716049396
$student.Cells.Item(11, 6).Value = 0.21353065539112046
$student.Cells.Item(12, 1).Value = "'11"
$student.Cells.Item(12, 2).Value = "Eta2"
$student.Cells.Item(12, 3).Value = 0.24580152671755728
$student.Cells.Item(12, 4).Value = 0.13283018867924531
$student.Cells.Item(12, 5).Value = 0.06510416666666673
$student.Cells.Item(12, 6).Value = 0.16365688487584648
$student.Cells.Item(13, 1).Value = "'12"
$student.Cells.Item(13, 2).Value = "DEL"
$student.Cells.Item(13, 3).Value = 0.00030364372469632284
$student.Cells.Item(13, 4).Value = -0.010232558139534904
$student.Cells.Item(13, 5).Value = -0.013215400624349567
$student.Cells.Item(13, 6).Value = -0.04394211471427059
$student.Cells.Item(14, 1).Value = "'13"
$student.Cells.Item(14, 2).Value = "ACO"
$student.Cells.Item(14, 3).Value = 0.0
$student.Cells.Item(14, 4).Value = 0.0013098236775818315
$student.Cells.Item(14, 5).Value = -0.0007014730934963763
$student.Cells.Item(14, 6).Value = -0.004218137993371479
$student.Cells.Item(15, 1).Value = "'14"
$student.Cells.Item(15, 2).Value = "ACL"
$student.Cells.Item(15, 3).Value = 0.39292364990689027
$student.Cells.Item(15, 4).Value = 0.25823045267489725
$student.Cells.Item(15, 5).Value = 0.1164215686274509
$student.Cells.Item(15, 6).Value = 0.36988543371522087
$student.Cells.Item(16, 1).Value = "'15"
$student.Cells.Item(16, 2).Value = "Pxx"
$student.Cells.Item(16, 3).Value = 0.010991102440881143
$student.Cells.Item(16, 4).Value = 0.032104492187499986
$student.Cells.Item(16, 5).Value = 0.1719641779963807
$student.Cells.Item(16, 6).Value = 0.2704822975646361
$student.Cells.Item(17, 1).Value = "'16"
$student.Cells.Item(17, 2).Value = "Pxx Exp(-Dij)"
$student.Cells.Item(17, 3).Value = -0.09491114701130854
$student.Cells.Item(17, 4).Value = 0.007417873542917664
$student.Cells.Item(17, 5).Value = -0.14880144730890993
$student.Cells.Item(17, 6).Value = -0.25226757369614516
$student.Cells.Item(18, 1).Value = "'17"
$student.Cells.Item(18, 2).Value = "DPxx"
$student.Cells.Item(18, 3).Value = 0.38888888888888884
$student.Cells.Item(18, 4).Value = 0.3008130081300813
$student.Cells.Item(18, 5).Value = 0.19540229885057475
$student.Cells.Item(18, 6).Value = 0.36038186157517904
$student.Cells.Item(19, 1).Value = "'18"
$student.Cells.Item(19, 2).Value = "PCC"
$student.Cells.Item(19, 3).Value = -0.2626262626262627
$student.Cells.Item(19, 4).Value = 0.061312607944732256
$student.Cells.Item(19, 5).Value = 1.7567567567567568
$student.Cells.Item(19, 6).Value = 2.9787234042553195
$student.Cells.Item(20, 1).Value = "'19"
$student.Cells.Item(20, 2).Value = "ACE"
$student.Cells.Item(20, 3).Value = 0.0019120458891013512
$student.Cells.Item(20, 4).Value = -0.002562788313685349
$student.Cells.Item(20, 5).Value = -0.004596527068437128
$student.Cells.Item(20, 6).Value = -0.015537332757876532
